# Update analysis with KM curve interpretation
#
# Inserts a new narrative paragraph (styled "First Paragraph", matching the
# style previously used by the paragraph that held the Kaplan-Meier figure)
# immediately before the Kaplan-Meier Curves figure, and re-styles the
# figure's paragraph to "Body Text" - mirroring the pattern used elsewhere
# in the document where a "First Paragraph" is followed by "Body Text"
# paragraphs.

$d = $word.ActiveDocument

# Locate the (single) paragraph that carries the inline KM chart - this is
# more robust than a hard-coded paragraph index.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the Kaplan-Meier figure paragraph"
}

$figurePara = $d.Paragraphs.Item($targetIndex)

# Insert a new empty paragraph right before the figure; it inherits the
# figure paragraph's current style ("First Paragraph").
$figurePara.Range.InsertParagraphBefore() | Out-Null

# The newly inserted paragraph is now at $targetIndex; the figure paragraph
# shifted down to $targetIndex + 1.
$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "During a median follow-up of 4.34 years [IQR: 2.29, 6.69], 19 deaths (21.6%) were observed among patients diagnosed with RCC. Of these, 5 deaths (12.8%) occurred in the pre-transplant group and 14 deaths (28.6%) in the post-transplant group. Kaplan-Meier survival analysis indicated a significantly higher risk of death in post-transplant RCC patients (Log-rank p = 0.047)."
$newPara.Style = "First Paragraph"

$figurePara = $d.Paragraphs.Item($targetIndex + 1)
$figurePara.Style = "Body Text"

Write-Output "Inserted KM interpretation paragraph before figure (index $targetIndex); figure paragraph restyled to Body Text."
